$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 43.231204986572266
$ws.Range("C2").Value = 6.632184028625488
$ws.Range("D2").Value = 27.603384017944336
$ws.Range("E2").Value = 57.85714340209961
